$d = $word.ActiveDocument

$replacements = @(
    @{old = "812×5=4060"; new = "978×9=8802"},
    @{old = "454×6=2724"; new = "725×6=4350"},
    @{old = "344×2=688";  new = "220×2=440"},
    @{old = "126×5=630";  new = "558×3=1674"},
    @{old = "285×3=855";  new = "646×4=2584"},
    @{old = "414×9=3726"; new = "352×6=2112"},
    @{old = "377×7=2639"; new = "818×6=4908"},
    @{old = "845×9=7605"; new = "121×8=968"},
    @{old = "823×5=4115"; new = "722×7=5054"},
    @{old = "638×8=5104"; new = "936×3=2808"},
    @{old = "544×3=1632"; new = "671×4=2684"},
    @{old = "334×7=2338"; new = "889×5=4445"},
    @{old = "788×4=3152"; new = "489×7=3423"},
    @{old = "585×8=4680"; new = "834×6=5004"},
    @{old = "908×6=5448"; new = "408×9=3672"},
    @{old = "809×2=1618"; new = "907×6=5442"},
    @{old = "334×3=1002"; new = "912×6=5472"},
    @{old = "348×9=3132"; new = "860×7=6020"},
    @{old = "683×2=1366"; new = "140×6=840"},
    @{old = "208×3=624";  new = "873×3=2619"},
    @{old = "898×6=5388"; new = "991×4=3964"},
    @{old = "353×4=1412"; new = "482×2=964"},
    @{old = "945×8=7560"; new = "591×9=5319"},
    @{old = "126×6=756";  new = "718×3=2154"},
    @{old = "197×9=1773"; new = "226×7=1582"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
